$d = $word.ActiveDocument
$d.Content.Find.Execute("εκτός έδρας για εκτέλεση υπηρεσίας", $true, $false, $false, $false, $false,
                         $true, 1, $false, "εντός έδρας για εκτέλεση υπηρεσίας", 2)
